$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), pushing every existing
# row down by one (the former last row, 159, becomes row 160).
$ws.Rows.Item(2).Insert()

# Populate the new top row with the latest price entry. The leading
# apostrophe keeps the date-shaped strings stored as text (matching the
# rest of the Date / Circular Date columns) instead of being parsed into
# date serials.
$ws.Range("A2").Value = "'12-01-2026"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 307.25
$ws.Range("E2").Value = "'01-01-2026"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

# Copy the formatting from the row below (the old row 2, now row 3) onto
# the newly inserted row so it matches the rest of the table exactly
# (this also clears the stray "quote prefix" flag picked up above).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The row that used to be the last row (159) is now row 160; it keeps its
# own text/values but the row-insert doesn't carry its hyperlink
# relationship down, so add it explicitly.
$ws.Hyperlinks.Add($ws.Range("F160"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Hyperlinks.Add applies Excel's default blue/underlined "Hyperlink" style;
# restore the plain look used by every other link cell in this sheet.
$ws.Range("F159").Copy()
$ws.Range("F160").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F160").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
